$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.137263894081116
$ws.Range("B1").Value = 1.568619132041931
$ws.Range("C1").Value = 2.720923662185669
$ws.Range("D1").Value = 4.178527355194092
$ws.Range("E1").Value = 1.681321740150452
